# "nuevos horarios rebrand AMC"
# Apply the new AMC-rebrand schedule values to the Data sheet, preserving
# the order in which new labels are first introduced so the shared-string
# table matches the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Row 3 --------------------------------------------------------------
$ws.Range("B3").Value = "LA VIDA OLALÁ SIN BUMP"

# --- Row 4 --------------------------------------------------------------
$ws.Range("B4").Value = "walking dead again and again con bump"

# --- Row 5 --------------------------------------------------------------
$ws.Range("B5").Value = "walking dead again and again sin bump"

# --- Row 6 --------------------------------------------------------------
$ws.Range("B6").Value = "GUARDIANES DE TRADICION con bump"
$ws.Range("A6").Value = "ee"

# --- Row 7 --------------------------------------------------------------
$ws.Range("B7").Value = "GUARDIANES DE TRADICION sin bump"
$ws.Range("A7").Value = "ee"

# --- Row 3 (feed column) -------------------------------------------------
$ws.Range("A3").Value = "mclatam"

# --- Row 4 (feed column) -------------------------------------------------
$ws.Range("A4").Value = "faLATAM"

# --- Row 5 (feed column) -------------------------------------------------
$ws.Range("A5").Value = "falatAM"

# --- Row 2 ----------------------------------------------------------------
$ws.Range("B2").Value = "fucking walking dead"
$ws.Range("A2").Value = "AMC*4FEEDS"

# --- Remaining column C/D/E value updates (reuse existing labels) --------
$ws.Range("C2").Value = "ESTRENO"
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 44116.979166666664

$ws.Range("C3").Value = "GEN"

$ws.Range("C4").Value = "PUNTUAL"

$ws.Range("C5").Value = "NT"

$ws.Range("C6").Value = "PELI DEL MES"

# --- Update the active selection to match the target workbook ------------
$ws.Activate() | Out-Null
$ws.Range("E3").Select() | Out-Null
